$wb = $excel.ActiveWorkbook
$wsGainers = $wb.Worksheets.Item("Top Gainers")
$wsLosers = $wb.Worksheets.Item("Top Losers")

# --- Update "Top Gainers" sheet: rows 2-75 (row 76 unchanged) ---
$wsGainers.Cells.Item(2, 2).Value = "SOLARWORLD"
$wsGainers.Cells.Item(2, 3).Value = 14.7722
$wsGainers.Cells.Item(2, 4).Value = 10.7603
$wsGainers.Cells.Item(2, 5).Value = 6.2751
$wsGainers.Cells.Item(3, 2).Value = "BLUEDART"
$wsGainers.Cells.Item(3, 3).Value = 14.447
$wsGainers.Cells.Item(3, 4).Value = 13.8712
$wsGainers.Cells.Item(3, 5).Value = 11.1257
$wsGainers.Cells.Item(4, 2).Value = "DREDGECORP"
$wsGainers.Cells.Item(4, 3).Value = 11.6405
$wsGainers.Cells.Item(4, 4).Value = 16.1254
$wsGainers.Cells.Item(4, 5).Value = 16.8846
$wsGainers.Cells.Item(5, 2).Value = "ADANIGREEN"
$wsGainers.Cells.Item(5, 3).Value = 11.4419
$wsGainers.Cells.Item(5, 4).Value = 8.534599999999999
$wsGainers.Cells.Item(5, 5).Value = 8.9838
$wsGainers.Cells.Item(6, 2).Value = "INFOBEAN"
$wsGainers.Cells.Item(6, 3).Value = 9.9924
$wsGainers.Cells.Item(6, 4).Value = 23.0892
$wsGainers.Cells.Item(6, 5).Value = 38.1616
$wsGainers.Cells.Item(7, 2).Value = "FIVESTAR"
$wsGainers.Cells.Item(7, 3).Value = 9.366
$wsGainers.Cells.Item(7, 4).Value = 9.376200000000001
$wsGainers.Cells.Item(7, 5).Value = 9.457700000000001
$wsGainers.Cells.Item(8, 2).Value = "VBL"
$wsGainers.Cells.Item(8, 3).Value = 9.347099999999999
$wsGainers.Cells.Item(8, 4).Value = 7.629
$wsGainers.Cells.Item(8, 5).Value = 11.9225
$wsGainers.Cells.Item(9, 2).Value = "GENUSPOWER"
$wsGainers.Cells.Item(9, 3).Value = 8.6149
$wsGainers.Cells.Item(9, 4).Value = 6.8639
$wsGainers.Cells.Item(9, 5).Value = 3.693
$wsGainers.Cells.Item(10, 2).Value = "HEG"
$wsGainers.Cells.Item(10, 3).Value = 7.601
$wsGainers.Cells.Item(10, 4).Value = 11.5862
$wsGainers.Cells.Item(10, 5).Value = 13.9731
$wsGainers.Cells.Item(11, 2).Value = "M&MFIN"
$wsGainers.Cells.Item(11, 3).Value = 7.5358
$wsGainers.Cells.Item(11, 4).Value = 8.167
$wsGainers.Cells.Item(11, 5).Value = 17.1449
$wsGainers.Cells.Item(12, 2).Value = "RPOWER"
$wsGainers.Cells.Item(12, 3).Value = 7.4405
$wsGainers.Cells.Item(12, 4).Value = 3.8274
$wsGainers.Cells.Item(12, 5).Value = 5.7935
$wsGainers.Cells.Item(13, 2).Value = "BUTTERFLY"
$wsGainers.Cells.Item(13, 3).Value = 7.2873
$wsGainers.Cells.Item(13, 4).Value = 10.1884
$wsGainers.Cells.Item(13, 5).Value = 12.7793
$wsGainers.Cells.Item(14, 2).Value = "WALCHANNAG"
$wsGainers.Cells.Item(14, 3).Value = 7.116
$wsGainers.Cells.Item(14, 4).Value = 4.6422
$wsGainers.Cells.Item(14, 5).Value = -4.5991
$wsGainers.Cells.Item(15, 2).Value = "TMB"
$wsGainers.Cells.Item(15, 3).Value = 7.0992
$wsGainers.Cells.Item(15, 4).Value = 10.9297
$wsGainers.Cells.Item(15, 5).Value = 18.3434
$wsGainers.Cells.Item(16, 2).Value = "UTKARSHBNK"
$wsGainers.Cells.Item(16, 3).Value = 6.5708
$wsGainers.Cells.Item(16, 4).Value = -4.3759
$wsGainers.Cells.Item(16, 5).Value = -1.0486
$wsGainers.Cells.Item(17, 2).Value = "SAIL"
$wsGainers.Cells.Item(17, 3).Value = 6.3559
$wsGainers.Cells.Item(17, 4).Value = 8.5741
$wsGainers.Cells.Item(17, 5).Value = 4.5211
$wsGainers.Cells.Item(18, 2).Value = "POKARNA"
$wsGainers.Cells.Item(18, 3).Value = 6.3061
$wsGainers.Cells.Item(18, 4).Value = -1.148
$wsGainers.Cells.Item(18, 5).Value = 19.1438
$wsGainers.Cells.Item(19, 2).Value = "FISCHER"
$wsGainers.Cells.Item(19, 3).Value = 6.1423
$wsGainers.Cells.Item(19, 4).Value = 11.1452
$wsGainers.Cells.Item(19, 5).Value = 4.3063
$wsGainers.Cells.Item(20, 2).Value = "SANDUMA"
$wsGainers.Cells.Item(20, 3).Value = 6.1288
$wsGainers.Cells.Item(20, 4).Value = 3.6403
$wsGainers.Cells.Item(20, 5).Value = 32.1943
$wsGainers.Cells.Item(21, 2).Value = "PDSL"
$wsGainers.Cells.Item(21, 3).Value = 6.0643
$wsGainers.Cells.Item(21, 4).Value = 4.0607
$wsGainers.Cells.Item(21, 5).Value = 9.9481
$wsGainers.Cells.Item(22, 2).Value = "VAIBHAVGBL"
$wsGainers.Cells.Item(22, 3).Value = 6.0228
$wsGainers.Cells.Item(22, 4).Value = 6.5513
$wsGainers.Cells.Item(22, 5).Value = 12.9078
$wsGainers.Cells.Item(23, 2).Value = "ABREL"
$wsGainers.Cells.Item(23, 3).Value = 5.7949
$wsGainers.Cells.Item(23, 4).Value = 6.5621
$wsGainers.Cells.Item(23, 5).Value = 6.1205
$wsGainers.Cells.Item(24, 2).Value = "JISLJALEQS"
$wsGainers.Cells.Item(24, 3).Value = 5.673
$wsGainers.Cells.Item(24, 4).Value = 4.9687
$wsGainers.Cells.Item(24, 5).Value = -1.1037
$wsGainers.Cells.Item(25, 2).Value = "IOC"
$wsGainers.Cells.Item(25, 3).Value = 5.4297
$wsGainers.Cells.Item(25, 4).Value = 8.339399999999999
$wsGainers.Cells.Item(25, 5).Value = 8.758900000000001
$wsGainers.Cells.Item(26, 2).Value = "GRAPHITE"
$wsGainers.Cells.Item(26, 3).Value = 5.3109
$wsGainers.Cells.Item(26, 4).Value = 11.7171
$wsGainers.Cells.Item(26, 5).Value = 11.9183
$wsGainers.Cells.Item(27, 2).Value = "VINCOFE"
$wsGainers.Cells.Item(27, 3).Value = 5.1427
$wsGainers.Cells.Item(27, 4).Value = 12.0998
$wsGainers.Cells.Item(27, 5).Value = 10.4529
$wsGainers.Cells.Item(28, 2).Value = "ADANIENSOL"
$wsGainers.Cells.Item(28, 3).Value = 5.1346
$wsGainers.Cells.Item(28, 4).Value = 2.5682
$wsGainers.Cells.Item(28, 5).Value = 11.0347
$wsGainers.Cells.Item(29, 2).Value = "EPACKPEB"
$wsGainers.Cells.Item(29, 3).Value = 5.0562
$wsGainers.Cells.Item(29, 4).Value = -1.4372
$wsGainers.Cells.Item(29, 5).Value = "N/A"
$wsGainers.Cells.Item(30, 2).Value = "SURYAROSNI"
$wsGainers.Cells.Item(30, 3).Value = 5.0009
$wsGainers.Cells.Item(30, 4).Value = 11.4048
$wsGainers.Cells.Item(30, 5).Value = 3.0388
$wsGainers.Cells.Item(31, 2).Value = "MEGASOFT"
$wsGainers.Cells.Item(31, 3).Value = 4.9974
$wsGainers.Cells.Item(31, 4).Value = 15.7588
$wsGainers.Cells.Item(31, 5).Value = 33.5271
$wsGainers.Cells.Item(32, 2).Value = "PROZONER"
$wsGainers.Cells.Item(32, 3).Value = 4.9921
$wsGainers.Cells.Item(32, 4).Value = 15.7468
$wsGainers.Cells.Item(32, 5).Value = 36.095
$wsGainers.Cells.Item(33, 2).Value = "STALLION"
$wsGainers.Cells.Item(33, 3).Value = 4.9914
$wsGainers.Cells.Item(33, 4).Value = -5.2229
$wsGainers.Cells.Item(33, 5).Value = 21.4391
$wsGainers.Cells.Item(34, 2).Value = "INDOTHAI"
$wsGainers.Cells.Item(34, 3).Value = 4.9883
$wsGainers.Cells.Item(34, 4).Value = 4.7163
$wsGainers.Cells.Item(34, 5).Value = 43.9974
$wsGainers.Cells.Item(35, 2).Value = "ABDL"
$wsGainers.Cells.Item(35, 3).Value = 4.821
$wsGainers.Cells.Item(35, 4).Value = 3.7277
$wsGainers.Cells.Item(35, 5).Value = 26.3345
$wsGainers.Cells.Item(36, 2).Value = "CELLO"
$wsGainers.Cells.Item(36, 3).Value = 4.7856
$wsGainers.Cells.Item(36, 4).Value = 3.6309
$wsGainers.Cells.Item(36, 5).Value = 13.4722
$wsGainers.Cells.Item(37, 2).Value = "STLTECH"
$wsGainers.Cells.Item(37, 3).Value = 4.5953
$wsGainers.Cells.Item(37, 4).Value = 1.5905
$wsGainers.Cells.Item(37, 5).Value = 7.7073
$wsGainers.Cells.Item(38, 2).Value = "BAJAJINDEF"
$wsGainers.Cells.Item(38, 3).Value = 4.565
$wsGainers.Cells.Item(38, 4).Value = 3.4852
$wsGainers.Cells.Item(38, 5).Value = 10.5031
$wsGainers.Cells.Item(39, 2).Value = "DATAMATICS"
$wsGainers.Cells.Item(39, 3).Value = 4.5435
$wsGainers.Cells.Item(39, 4).Value = 6.95
$wsGainers.Cells.Item(39, 5).Value = 15.336
$wsGainers.Cells.Item(40, 2).Value = "RELINFRA"
$wsGainers.Cells.Item(40, 3).Value = 4.4466
$wsGainers.Cells.Item(40, 4).Value = -2.7162
$wsGainers.Cells.Item(40, 5).Value = -7.5382
$wsGainers.Cells.Item(41, 2).Value = "GMBREW"
$wsGainers.Cells.Item(41, 3).Value = 4.4121
$wsGainers.Cells.Item(41, 4).Value = -0.0396
$wsGainers.Cells.Item(41, 5).Value = 79.9117
$wsGainers.Cells.Item(42, 2).Value = "ORIENTTECH"
$wsGainers.Cells.Item(42, 3).Value = 4.3801
$wsGainers.Cells.Item(42, 4).Value = 1.0602
$wsGainers.Cells.Item(42, 5).Value = 33.3852
$wsGainers.Cells.Item(43, 2).Value = "JKIL"
$wsGainers.Cells.Item(43, 3).Value = 4.2818
$wsGainers.Cells.Item(43, 4).Value = 3.0893
$wsGainers.Cells.Item(43, 5).Value = 1.8997
$wsGainers.Cells.Item(44, 2).Value = "AXISCADES"
$wsGainers.Cells.Item(44, 3).Value = 4.2661
$wsGainers.Cells.Item(44, 4).Value = 6.7266
$wsGainers.Cells.Item(44, 5).Value = -3.2497
$wsGainers.Cells.Item(45, 2).Value = "PVRINOX"
$wsGainers.Cells.Item(45, 3).Value = 4.2469
$wsGainers.Cells.Item(45, 4).Value = 6.348
$wsGainers.Cells.Item(45, 5).Value = 14.8558
$wsGainers.Cells.Item(46, 2).Value = "PROSTARM"
$wsGainers.Cells.Item(46, 3).Value = 4.2353
$wsGainers.Cells.Item(46, 4).Value = 1.3441
$wsGainers.Cells.Item(46, 5).Value = -7.6164
$wsGainers.Cells.Item(47, 2).Value = "SGMART"
$wsGainers.Cells.Item(47, 3).Value = 4.2313
$wsGainers.Cells.Item(47, 4).Value = 8.215
$wsGainers.Cells.Item(47, 5).Value = 2.4965
$wsGainers.Cells.Item(48, 2).Value = "GPPL"
$wsGainers.Cells.Item(48, 3).Value = 4.1889
$wsGainers.Cells.Item(48, 4).Value = 3.183
$wsGainers.Cells.Item(48, 5).Value = 4.8219
$wsGainers.Cells.Item(49, 2).Value = "ATGL"
$wsGainers.Cells.Item(49, 3).Value = 4.0258
$wsGainers.Cells.Item(49, 4).Value = 3.7668
$wsGainers.Cells.Item(49, 5).Value = 3.3104
$wsGainers.Cells.Item(50, 2).Value = "FILATEX"
$wsGainers.Cells.Item(50, 3).Value = 4.0221
$wsGainers.Cells.Item(50, 4).Value = 9.383599999999999
$wsGainers.Cells.Item(50, 5).Value = 24.9853
$wsGainers.Cells.Item(51, 2).Value = "GPIL"
$wsGainers.Cells.Item(51, 3).Value = 4.0079
$wsGainers.Cells.Item(51, 4).Value = 6.175
$wsGainers.Cells.Item(51, 5).Value = 14.277
$wsGainers.Cells.Item(52, 2).Value = "HITECHGEAR"
$wsGainers.Cells.Item(52, 3).Value = 3.9677
$wsGainers.Cells.Item(52, 4).Value = 1.2548
$wsGainers.Cells.Item(52, 5).Value = 10.0407
$wsGainers.Cells.Item(53, 2).Value = "TCI"
$wsGainers.Cells.Item(53, 3).Value = 3.9548
$wsGainers.Cells.Item(53, 4).Value = 3.8585
$wsGainers.Cells.Item(53, 5).Value = 4.3596
$wsGainers.Cells.Item(54, 2).Value = "SRM"
$wsGainers.Cells.Item(54, 3).Value = 3.9497
$wsGainers.Cells.Item(54, 4).Value = 3.6408
$wsGainers.Cells.Item(54, 5).Value = 4.5535
$wsGainers.Cells.Item(55, 2).Value = "MRPL"
$wsGainers.Cells.Item(55, 3).Value = 3.9241
$wsGainers.Cells.Item(55, 4).Value = 9.352399999999999
$wsGainers.Cells.Item(55, 5).Value = 19.6626
$wsGainers.Cells.Item(56, 2).Value = "RHIM"
$wsGainers.Cells.Item(56, 3).Value = 3.9115
$wsGainers.Cells.Item(56, 4).Value = 3.4836
$wsGainers.Cells.Item(56, 5).Value = 5.4434
$wsGainers.Cells.Item(57, 2).Value = "CMSINFO"
$wsGainers.Cells.Item(57, 3).Value = 3.8961
$wsGainers.Cells.Item(57, 4).Value = 2.6738
$wsGainers.Cells.Item(57, 5).Value = 2.8801
$wsGainers.Cells.Item(58, 2).Value = "HCC"
$wsGainers.Cells.Item(58, 3).Value = 3.8873
$wsGainers.Cells.Item(58, 4).Value = 2.7876
$wsGainers.Cells.Item(58, 5).Value = 7.5305
$wsGainers.Cells.Item(59, 2).Value = "LLOYDSENT"
$wsGainers.Cells.Item(59, 3).Value = 3.8566
$wsGainers.Cells.Item(59, 4).Value = 1.1444
$wsGainers.Cells.Item(59, 5).Value = 10.4808
$wsGainers.Cells.Item(60, 2).Value = "RECLTD"
$wsGainers.Cells.Item(60, 3).Value = 3.8488
$wsGainers.Cells.Item(60, 4).Value = 2.8315
$wsGainers.Cells.Item(60, 5).Value = 2.7625
$wsGainers.Cells.Item(61, 2).Value = "NBCC"
$wsGainers.Cells.Item(61, 3).Value = 3.8165
$wsGainers.Cells.Item(61, 4).Value = 2.5338
$wsGainers.Cells.Item(61, 5).Value = 6.9481
$wsGainers.Cells.Item(62, 2).Value = "MAHLOG"
$wsGainers.Cells.Item(62, 3).Value = 3.7305
$wsGainers.Cells.Item(62, 4).Value = 0.6736
$wsGainers.Cells.Item(62, 5).Value = 4.2127
$wsGainers.Cells.Item(63, 2).Value = "INOXWIND"
$wsGainers.Cells.Item(63, 3).Value = 3.7247
$wsGainers.Cells.Item(63, 4).Value = 3.1996
$wsGainers.Cells.Item(63, 5).Value = 13.3842
$wsGainers.Cells.Item(64, 2).Value = "STAR"
$wsGainers.Cells.Item(64, 3).Value = 3.6855
$wsGainers.Cells.Item(64, 4).Value = 3.6155
$wsGainers.Cells.Item(64, 5).Value = 2.8516
$wsGainers.Cells.Item(65, 2).Value = "APARINDS"
$wsGainers.Cells.Item(65, 3).Value = 3.6783
$wsGainers.Cells.Item(65, 4).Value = 8.1182
$wsGainers.Cells.Item(65, 5).Value = 15.3494
$wsGainers.Cells.Item(66, 2).Value = "SUNFLAG"
$wsGainers.Cells.Item(66, 3).Value = 3.656
$wsGainers.Cells.Item(66, 4).Value = 3.9909
$wsGainers.Cells.Item(66, 5).Value = 4.2882
$wsGainers.Cells.Item(67, 2).Value = "SUZLON"
$wsGainers.Cells.Item(67, 3).Value = 3.6464
$wsGainers.Cells.Item(67, 4).Value = 8.2683
$wsGainers.Cells.Item(67, 5).Value = 5.8492
$wsGainers.Cells.Item(68, 2).Value = "ICRA"
$wsGainers.Cells.Item(68, 3).Value = 3.6236
$wsGainers.Cells.Item(68, 4).Value = 4.3033
$wsGainers.Cells.Item(68, 5).Value = 2.7095
$wsGainers.Cells.Item(69, 2).Value = "NMDC"
$wsGainers.Cells.Item(69, 3).Value = 3.6179
$wsGainers.Cells.Item(69, 4).Value = 4.2324
$wsGainers.Cells.Item(69, 5).Value = 1.2836
$wsGainers.Cells.Item(70, 2).Value = "SAMBHV"
$wsGainers.Cells.Item(70, 3).Value = 3.6169
$wsGainers.Cells.Item(70, 4).Value = 2.1136
$wsGainers.Cells.Item(70, 5).Value = 4.6439
$wsGainers.Cells.Item(71, 2).Value = "ASHAPURMIN"
$wsGainers.Cells.Item(71, 3).Value = 3.6026
$wsGainers.Cells.Item(71, 4).Value = 6.254
$wsGainers.Cells.Item(71, 5).Value = 2.043
$wsGainers.Cells.Item(72, 2).Value = "SHK"
$wsGainers.Cells.Item(72, 3).Value = 3.4745
$wsGainers.Cells.Item(72, 4).Value = 2.2297
$wsGainers.Cells.Item(72, 5).Value = -2.0836
$wsGainers.Cells.Item(73, 2).Value = "GAIL"
$wsGainers.Cells.Item(73, 3).Value = 3.4686
$wsGainers.Cells.Item(73, 4).Value = 2.0053
$wsGainers.Cells.Item(73, 5).Value = 4.7422
$wsGainers.Cells.Item(74, 2).Value = "IVALUE"
$wsGainers.Cells.Item(74, 3).Value = 3.4675
$wsGainers.Cells.Item(74, 4).Value = 3.2487
$wsGainers.Cells.Item(74, 5).Value = -3.7455
$wsGainers.Cells.Item(75, 2).Value = "MSTCLTD"
$wsGainers.Cells.Item(75, 3).Value = 3.4612
$wsGainers.Cells.Item(75, 4).Value = 3.3933
$wsGainers.Cells.Item(75, 5).Value = 15.7641

# --- Update "Top Losers" sheet ---
$wsLosers.Cells.Item(10, 4).Value = 5.978
$wsLosers.Cells.Item(73, 4).Value = "N/A"

Write-Host "Edit complete"